$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update DAMSLTag (column I) and DialogAct (column J) values for the rows
# that were re-annotated by re-running SGNN after transcript cleanup.
$ws.Range("I11").Value = "sv"
$ws.Range("J11").Value = "Statement-opinion"
$ws.Range("I27").Value = "aa"
$ws.Range("J27").Value = "Agree/Accept"
$ws.Range("I35").Value = "aa"
$ws.Range("J35").Value = "Agree/Accept"
$ws.Range("I42").Value = "b"
$ws.Range("J42").Value = "Acknowledge (Backchannel)"
$ws.Range("I57").Value = "sv"
$ws.Range("J57").Value = "Statement-opinion"
$ws.Range("I61").Value = "aa"
$ws.Range("J61").Value = "Agree/Accept"
$ws.Range("I69").Value = "b"
$ws.Range("J69").Value = "Acknowledge (Backchannel)"
$ws.Range("I79").Value = "sv"
$ws.Range("J79").Value = "Statement-opinion"
$ws.Range("I82").Value = "b"
$ws.Range("J82").Value = "Acknowledge (Backchannel)"
$ws.Range("I93").Value = "ba"
$ws.Range("J93").Value = "Appreciation"
$ws.Range("I96").Value = "ba"
$ws.Range("J96").Value = "Appreciation"
$ws.Range("I102").Value = "b"
$ws.Range("J102").Value = "Acknowledge (Backchannel)"
$ws.Range("I105").Value = "ba"
$ws.Range("J105").Value = "Appreciation"
$ws.Range("I113").Value = "sd"
$ws.Range("J113").Value = "Statement-non-opinion"
$ws.Range("I114").Value = "sd"
$ws.Range("J114").Value = "Statement-non-opinion"
$ws.Range("I117").Value = "aa"
$ws.Range("J117").Value = "Agree/Accept"
$ws.Range("I142").Value = "aa"
$ws.Range("J142").Value = "Agree/Accept"
$ws.Range("I144").Value = "ba"
$ws.Range("J144").Value = "Appreciation"
$ws.Range("I168").Value = "sd"
$ws.Range("J168").Value = "Statement-non-opinion"
$ws.Range("I170").Value = "sd"
$ws.Range("J170").Value = "Statement-non-opinion"
$ws.Range("I174").Value = "sd"
$ws.Range("J174").Value = "Statement-non-opinion"
$ws.Range("I179").Value = "ba"
$ws.Range("J179").Value = "Appreciation"
$ws.Range("I200").Value = "sd"
$ws.Range("J200").Value = "Statement-non-opinion"
$ws.Range("I209").Value = "sv"
$ws.Range("J209").Value = "Statement-opinion"
$ws.Range("I247").Value = "ba"
$ws.Range("J247").Value = "Appreciation"
$ws.Range("I261").Value = "aa"
$ws.Range("J261").Value = "Agree/Accept"
$ws.Range("I267").Value = "sv"
$ws.Range("J267").Value = "Statement-opinion"
$ws.Range("I278").Value = "sv"
$ws.Range("J278").Value = "Statement-opinion"
$ws.Range("I279").Value = "sd"
$ws.Range("J279").Value = "Statement-non-opinion"
$ws.Range("I290").Value = "sd"
$ws.Range("J290").Value = "Statement-non-opinion"
$ws.Range("I291").Value = "aa"
$ws.Range("J291").Value = "Agree/Accept"
$ws.Range("I295").Value = "aa"
$ws.Range("J295").Value = "Agree/Accept"
$ws.Range("I305").Value = "sd"
$ws.Range("J305").Value = "Statement-non-opinion"
$ws.Range("I308").Value = "%"
$ws.Range("J308").Value = "Uninterpretable"
$ws.Range("I316").Value = "sv"
$ws.Range("J316").Value = "Statement-opinion"
$ws.Range("I318").Value = "%"
$ws.Range("J318").Value = "Uninterpretable"
$ws.Range("I319").Value = "sd"
$ws.Range("J319").Value = "Statement-non-opinion"
$ws.Range("I322").Value = "%"
$ws.Range("J322").Value = "Uninterpretable"
$ws.Range("I324").Value = "sd"
$ws.Range("J324").Value = "Statement-non-opinion"
$ws.Range("I328").Value = "aa"
$ws.Range("J328").Value = "Agree/Accept"
$ws.Range("I336").Value = "sd"
$ws.Range("J336").Value = "Statement-non-opinion"
$ws.Range("I343").Value = "b"
$ws.Range("J343").Value = "Acknowledge (Backchannel)"
$ws.Range("I350").Value = "b"
$ws.Range("J350").Value = "Acknowledge (Backchannel)"
$ws.Range("I353").Value = "sv"
$ws.Range("J353").Value = "Statement-opinion"
$ws.Range("I366").Value = "ba"
$ws.Range("J366").Value = "Appreciation"
$ws.Range("I377").Value = "ba"
$ws.Range("J377").Value = "Appreciation"
$ws.Range("I382").Value = "sd"
$ws.Range("J382").Value = "Statement-non-opinion"
$ws.Range("I385").Value = "%"
$ws.Range("J385").Value = "Uninterpretable"
$ws.Range("I386").Value = "aa"
$ws.Range("J386").Value = "Agree/Accept"
$ws.Range("I387").Value = "sv"
$ws.Range("J387").Value = "Statement-opinion"
$ws.Range("I408").Value = "sv"
$ws.Range("J408").Value = "Statement-opinion"
$ws.Range("I433").Value = "aa"
$ws.Range("J433").Value = "Agree/Accept"
$ws.Range("I435").Value = "sd"
$ws.Range("J435").Value = "Statement-non-opinion"
$ws.Range("I445").Value = "b"
$ws.Range("J445").Value = "Acknowledge (Backchannel)"
$ws.Range("I446").Value = "sd"
$ws.Range("J446").Value = "Statement-non-opinion"
$ws.Range("I457").Value = "b"
$ws.Range("J457").Value = "Acknowledge (Backchannel)"
$ws.Range("I462").Value = "sd"
$ws.Range("J462").Value = "Statement-non-opinion"
$ws.Range("I463").Value = "sd"
$ws.Range("J463").Value = "Statement-non-opinion"
$ws.Range("I465").Value = "%"
$ws.Range("J465").Value = "Uninterpretable"
$ws.Range("I468").Value = "b"
$ws.Range("J468").Value = "Acknowledge (Backchannel)"
$ws.Range("I469").Value = "sd"
$ws.Range("J469").Value = "Statement-non-opinion"
$ws.Range("I481").Value = "sv"
$ws.Range("J481").Value = "Statement-opinion"
$ws.Range("I488").Value = "ba"
$ws.Range("J488").Value = "Appreciation"
